$wb = $excel.ActiveWorkbook

# --- Sheet: MAIN_CONTROLLER ---
$ws1 = $wb.Worksheets.Item("MAIN_CONTROLLER")
$ws1.Range("D1").Value = "ApplicationName"

# --- Sheet: DATASHEET ---
$ws2 = $wb.Worksheets.Item("DATASHEET")
$ws2.Range("C2").Value = "FOSs"
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Y"
$ws2.Range("C3").Value = "FOS"
$ws2.Range("D3").Value = "FOS.xlsx"
$ws2.Range("E3").Value = 20
$ws2.Range("F3").Value = 30

# Update selections to match the saved UI state (DATASHEET selection first,
# then re-activate MAIN_CONTROLLER so it remains the selected/visible tab).
$ws2.Activate()
$ws2.Range("B3").Select() | Out-Null

$ws1.Activate()
$ws1.Range("E5").Select() | Out-Null

$wb.Save()
